$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.86
$summary.Range("B4").Value = -0.14
$summary.Range("B5").Value = -1.4
$summary.Range("B6").Value = 2
$summary.Range("B8").Value = 2

# ---- Strategy Status sheet ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.86
$status.Range("D4").Value = 2
$status.Range("E4").Value = -0.14
$status.Range("F4").Value = -0.14

# ---- New trade row data (Trade #2) ----
function Add-TradeRow($ws) {
    $ws.Cells.Item(3, 1).Value = 2
    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 2).Value = "2026-02-17"
    $ws.Cells.Item(3, 2).Style = "Normal"
    $ws.Cells.Item(3, 3).Value = "13:01:56"
    $ws.Cells.Item(3, 4).Value = "MarketMaking"
    $ws.Cells.Item(3, 5).Value = "DOWN"
    $ws.Cells.Item(3, 6).Value = 0.59
    $ws.Cells.Item(3, 7).Value = 0.48
    $ws.Cells.Item(3, 8).Value = "CLOSED"
    $ws.Cells.Item(3, 9).Value = -18.6441
    $ws.Cells.Item(3, 10).Value = -0.11
    $ws.Cells.Item(3, 11).Value = 99.86
    $ws.Cells.Item(3, 12).Value = 0
    $ws.Cells.Item(3, 13).Value = 0
    $ws.Cells.Item(3, 14).Value = 0.6
    $ws.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(3, 16).Value = "early_exit"
    $ws.Cells.Item(3, 17).Value = 0.13
}

# ---- All Trades sheet ----
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# ---- MarketMaking sheet ----
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
